# Swap the deck's theme palette from "Integral" to the stock "Office Theme"
# (ppt/theme/theme1.xml), matching the authored commit that exchanged the
# colour schemes used by theme1.xml / theme2.xml.
#
# NOTE: fontScheme / fmtScheme are byte-identical between the two themes in
# this deck, so only the 12 clrScheme entries actually need to change.

$p = $ppt.ActivePresentation

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeTheme = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

$cs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $rgb = $officeTheme[$i]
    $cs.Item($i + 1).RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}
